$wb = $excel.ActiveWorkbook

# --- Sheet: normality ---
$ws2 = $wb.Worksheets.Item("normality")
$ws2.Range("C3").Value = 0.9873
$ws2.Range("D3").Value = 0.3637
$ws2.Range("C6").Value = 0.9875
$ws2.Range("D6").Value = 0.3781
$ws2.Range("C9").Value = 0.9556
$ws2.Range("D9").Value = 0.0008
$ws2.Range("C12").Value = 0.9648
$ws2.Range("D12").Value = 0.0043
$ws2.Range("D15").Value = 0.07149999999999999
$ws2.Range("C18").Value = 0.9712
$ws2.Range("D18").Value = 0.0147
$ws2.Range("C21").Value = 0.9681999999999999
$ws2.Range("D21").Value = 0.008200000000000001

# --- Sheet: equal_var ---
$ws3 = $wb.Worksheets.Item("equal_var")
$ws3.Range("C3").Value = 2.0444
$ws3.Range("D3").Value = 0.1327
$ws3.Range("C4").Value = 0.3262
$ws3.Range("D4").Value = 0.7222
$ws3.Range("C5").Value = 1.1509
$ws3.Range("D5").Value = 0.3189
$ws3.Range("C6").Value = 0.0171
$ws3.Range("D6").Value = 0.983
$ws3.Range("C7").Value = 1.1119
$ws3.Range("D7").Value = 0.3314
$ws3.Range("C8").Value = 0.1466
$ws3.Range("D8").Value = 0.8638
$ws3.Range("C9").Value = 0.4392
$ws3.Range("D9").Value = 0.6453

# --- Sheet: anova ---
$ws4 = $wb.Worksheets.Item("anova")
$ws4.Range("E3").Value = 166
$ws4.Range("F3").Value = 0.9429
$ws4.Range("G3").Value = 0.3916
$ws4.Range("H3").Value = 0.0112
$ws4.Range("E4").Value = 166
$ws4.Range("F4").Value = 2.0563
$ws4.Range("G4").Value = 0.1312
$ws4.Range("H4").Value = 0.0242
$ws4.Range("E5").Value = 166
$ws4.Range("F5").Value = 0.8853
$ws4.Range("G5").Value = 0.4145
$ws4.Range("H5").Value = 0.0106
$ws4.Range("E6").Value = 166
$ws4.Range("F6").Value = 2.3991
$ws4.Range("G6").Value = 0.0939
$ws4.Range("H6").Value = 0.0281
$ws4.Range("E7").Value = 166
$ws4.Range("F7").Value = 1.5656
$ws4.Range("G7").Value = 0.212
$ws4.Range("H7").Value = 0.0185
$ws4.Range("E8").Value = 166
$ws4.Range("F8").Value = 2.647
$ws4.Range("G8").Value = 0.07389999999999999
$ws4.Range("H8").Value = 0.0309
$ws4.Range("E9").Value = 166
$ws4.Range("F9").Value = 2.9356
$ws4.Range("G9").Value = 0.0559
$ws4.Range("H9").Value = 0.0342

# --- Sheet: pairwise_ttests ---
$ws5 = $wb.Worksheets.Item("pairwise_ttests")
$ws5.Range("H3").Value = -1.1943
$ws5.Range("I3").Value = 60.5547
$ws5.Range("K3").Value = 0.237
$ws5.Range("M3").Value = -0.2422
$ws5.Range("N3").Value = 0.711
$ws5.Range("H5").Value = 0.3995
$ws5.Range("I5").Value = 13.1532
$ws5.Range("K5").Value = 0.6959
$ws5.Range("M5").Value = 0.1262
$ws5.Range("H6").Value = -1.8607
$ws5.Range("I6").Value = 70.0136
$ws5.Range("K6").Value = 0.067
$ws5.Range("M6").Value = -0.3453
$ws5.Range("N6").Value = 0.201
$ws5.Range("H8").Value = -0.2352
$ws5.Range("I8").Value = 14.374
$ws5.Range("K8").Value = 0.8174
$ws5.Range("M8").Value = -0.0624
$ws5.Range("H9").Value = 0.8364
$ws5.Range("I9").Value = 73.1742
$ws5.Range("K9").Value = 0.4056
$ws5.Range("M9").Value = 0.1516
$ws5.Range("H11").Value = 1.156
$ws5.Range("I11").Value = 16.1144
$ws5.Range("K11").Value = 0.2645
$ws5.Range("M11").Value = 0.2628
$ws5.Range("N11").Value = 0.7935000000000001
$ws5.Range("H12").Value = -2.1348
$ws5.Range("I12").Value = 73.6656
$ws5.Range("K12").Value = 0.0361
$ws5.Range("M12").Value = -0.3855
$ws5.Range("N12").Value = 0.1083
$ws5.Range("H14").Value = 0.0426
$ws5.Range("I14").Value = 13.5152
$ws5.Range("K14").Value = 0.9666
$ws5.Range("M14").Value = 0.0127
$ws5.Range("H15").Value = -1.5537
$ws5.Range("I15").Value = 62.0384
$ws5.Range("K15").Value = 0.1254
$ws5.Range("M15").Value = -0.31
$ws5.Range("N15").Value = 0.3762
$ws5.Range("H17").Value = 0.599
$ws5.Range("I17").Value = 13.0163
$ws5.Range("K17").Value = 0.5594
$ws5.Range("M17").Value = 0.1942
$ws5.Range("H18").Value = -2.2211
$ws5.Range("I18").Value = 73.8891
$ws5.Range("K18").Value = 0.0294
$ws5.Range("M18").Value = -0.4004
$ws5.Range("N18").Value = 0.0882
$ws5.Range("H20").Value = -0.0669
$ws5.Range("I20").Value = 13.8322
$ws5.Range("K20").Value = 0.9476
$ws5.Range("M20").Value = -0.019
$ws5.Range("H21").Value = -2.2231
$ws5.Range("I21").Value = 73.48950000000001
$ws5.Range("K21").Value = 0.0293
$ws5.Range("M21").Value = -0.4019
$ws5.Range("N21").Value = 0.08790000000000001
$ws5.Range("H23").Value = -0.5052
$ws5.Range("I23").Value = 12.8451
$ws5.Range("K23").Value = 0.622
$ws5.Range("M23").Value = -0.1696

# BF10 column (L) holds text-typed numeric-looking strings; force text format
# so Excel does not auto-convert them to numbers, preserving original string cell type.
$ws5.Range("L3").NumberFormat = "@"
$ws5.Range("L3").Value = "0.364"
$ws5.Range("L4").NumberFormat = "@"
$ws5.Range("L4").Value = "0.335"
$ws5.Range("L5").NumberFormat = "@"
$ws5.Range("L5").Value = "0.319"
$ws5.Range("L6").NumberFormat = "@"
$ws5.Range("L6").Value = "0.909"
$ws5.Range("L7").NumberFormat = "@"
$ws5.Range("L7").Value = "0.683"
$ws5.Range("L8").NumberFormat = "@"
$ws5.Range("L8").Value = "0.306"
$ws5.Range("L9").NumberFormat = "@"
$ws5.Range("L9").Value = "0.262"
$ws5.Range("L10").NumberFormat = "@"
$ws5.Range("L10").Value = "0.858"
$ws5.Range("L11").NumberFormat = "@"
$ws5.Range("L11").Value = "0.51"
$ws5.Range("L12").NumberFormat = "@"
$ws5.Range("L12").Value = "1.483"
$ws5.Range("L13").NumberFormat = "@"
$ws5.Range("L13").Value = "0.538"
$ws5.Range("L14").NumberFormat = "@"
$ws5.Range("L14").Value = "0.299"
$ws5.Range("L15").NumberFormat = "@"
$ws5.Range("L15").Value = "0.568"
$ws5.Range("L16").NumberFormat = "@"
$ws5.Range("L16").Value = "0.335"
$ws5.Range("L17").NumberFormat = "@"
$ws5.Range("L17").Value = "0.345"
$ws5.Range("L18").NumberFormat = "@"
$ws5.Range("L18").Value = "1.753"
$ws5.Range("L19").NumberFormat = "@"
$ws5.Range("L19").Value = "0.656"
$ws5.Range("L20").NumberFormat = "@"
$ws5.Range("L20").Value = "0.3"
$ws5.Range("L21").NumberFormat = "@"
$ws5.Range("L21").Value = "1.76"
$ws5.Range("L22").NumberFormat = "@"
$ws5.Range("L22").Value = "0.847"
$ws5.Range("L23").NumberFormat = "@"
$ws5.Range("L23").Value = "0.331"

Write-Output "Edit complete"
